$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set F3 to "Fail" (Pass/Fail column result for the failing test case)
$ws.Range("F3").Value = "Fail"

# Update G3 message from "-" to "Hibás"
$ws.Range("G3").Value = "Hibás"

# Update the selection on the sheet to I10
$ws.Range("I10").Select()
